# asImageInHeader-template.docx — M2Doc issue #295
#
# The canonical-OOXML diff for this particular template shows no textual,
# formatting, or structural change at all: every "-"/"+" pair in
# word/document.xml, word/footnotes.xml, word/header1.xml and
# word/styles.xml is the very same element with the very same attribute
# name/value pairs, only re-ordered (e.g. <w:color w:val="E36C0A"
# w:themeColor="accent6" w:themeShade="BF"/> becomes <w:color
# w:themeColor="accent6" w:themeShade="BF" w:val="E36C0A"/>, the root
# namespace declarations get alphabetised, rsid-style bookkeeping
# attributes are untouched, etc.). That reordering is a side effect of the
# tooling that regenerated the test-resource templates while the fix for
# #295 (stamping the M2Doc version into a template's custom document
# properties) was implemented — it is not a content edit of this
# particular fixture, which has no custom-properties part before or after
# the commit.
#
# So there is nothing in this template's body/header/styles for a Word
# automation script to legitimately rewrite — doing so (e.g. touching the
# header's field-code runs, the bookmark, or the style catalog) would
# introduce changes that are not actually present in the diff. The
# correct COM-interop action is therefore to just touch the document
# through the object model (proving the session is live / the expected
# parts are present) without mutating any text, formatting or structure.

$d = $word.ActiveDocument

# Sanity: confirm the parts this commit's diff concerns are still the
# ones we expect — one paragraph, one section, one header, the "_GoBack"
# bookmark — then leave the content exactly as-is.
$paragraphCount = $d.Paragraphs.Count
$sectionCount = $d.Sections.Count
$section = $d.Sections.Item(1)
$header = $section.Headers.Item(1)
$hasGoBack = $d.Bookmarks.Exists("_GoBack")

Write-Output "paragraphs=$paragraphCount sections=$sectionCount header=$($header -ne $null) _GoBack=$hasGoBack"
